# [Kadastro App] Yeni kayit eklendi: 2907
$wb = $excel.ActiveWorkbook

# New record data, appended identically to the "Kayitlar" summary sheet
# and the "Erdemli" unit sheet.
$recordNo   = "2907"
$recordDate = "2025-09-08"
$recordUnit = "Erdemli"
$parselSayi = "1"
$recordJob  = "3B"
$personnel  = "EMİNE ALANLI KIRCILI (K.Mühendisi), SEVİL SARAÇER (Tekniker)"

$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$rowKayitlar = 14
$wsKayitlar.Range("A" + $rowKayitlar + ":F" + $rowKayitlar).NumberFormat = "@"
$wsKayitlar.Cells.Item($rowKayitlar, 1).Value = $recordNo
$wsKayitlar.Cells.Item($rowKayitlar, 2).Value = $recordDate
$wsKayitlar.Cells.Item($rowKayitlar, 3).Value = $recordUnit
$wsKayitlar.Cells.Item($rowKayitlar, 4).Value = $parselSayi
$wsKayitlar.Cells.Item($rowKayitlar, 5).Value = $recordJob
$wsKayitlar.Cells.Item($rowKayitlar, 6).Value = $personnel

$wsErdemli = $wb.Worksheets.Item("Erdemli")
$rowErdemli = 13
$wsErdemli.Range("A" + $rowErdemli + ":F" + $rowErdemli).NumberFormat = "@"
$wsErdemli.Cells.Item($rowErdemli, 1).Value = $recordNo
$wsErdemli.Cells.Item($rowErdemli, 2).Value = $recordDate
$wsErdemli.Cells.Item($rowErdemli, 3).Value = $recordUnit
$wsErdemli.Cells.Item($rowErdemli, 4).Value = $parselSayi
$wsErdemli.Cells.Item($rowErdemli, 5).Value = $recordJob
$wsErdemli.Cells.Item($rowErdemli, 6).Value = $personnel

Write-Output "Added record $recordNo to Kayitlar!A$rowKayitlar and Erdemli!A$rowErdemli"
